$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.464.38'
$ws.Range('E2').Value = '  -0.69%  '
$ws.Range('D3').Value = '1.825.13'
$ws.Range('E3').Value = '  -1.00%  '
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').Value = '312.74'
$ws.Range('D6').Value = '1.001'
$ws.Range('E6').Value = '  -0.09%  '
$ws.Range('D7').Value = '0.4246'
$ws.Range('E7').Value = '  -0.45%  '
$ws.Range('D8').Value = '0.3618'
$ws.Range('E8').Value = '  +0.06%  '
$ws.Range('D9').Value = '0.07213'
$ws.Range('E9').Value = '  -1.28%  '
$ws.Range('D10').Value = '0.8605'
$ws.Range('E10').Value = '  -1.06%  '
$ws.Range('D11').Value = '20.61'
$ws.Range('E11').Value = '  -0.35%  '
$ws.Range('D12').Value = '1.855.05'
$ws.Range('E12').Value = '  -0.24%  '
$ws.Range('D13').Value = '5.393'
$ws.Range('E13').Value = '  +1.22%  '
$ws.Range('D14').Value = '6.480'
$ws.Range('E14').Value = '  -0.87%  '
$ws.Range('D15').Value = '0.06927'
$ws.Range('E15').Value = '  -0.87%  '
$ws.Range('D16').Value = '1.003'
$ws.Range('E16').Value = '  -0.10%  '
$ws.Range('D17').Value = '80.35'
$ws.Range('E17').Value = '  +1.22%  '
$ws.Range('D18').Value = '0.000008882'
$ws.Range('E18').Value = '  -0.89%  '
$ws.Range('E19').Value = '  -0.03%  '
$ws.Range('D20').Value = '15.34'
$ws.Range('E20').Value = '  +0.36%  '
$ws.Range('D21').Value = '27.487.64'
$ws.Range('E21').Value = '  -0.76%  '
$ws.Range('D22').Value = '5.130'
$ws.Range('E22').Value = '  +3.03%  '
$ws.Range('D23').Value = '10.96'
$ws.Range('E23').Value = '  +5.94%  '
$ws.Range('D24').Value = '2.057.89'
$ws.Range('E24').Value = '  -1.33%  '
$ws.Range('E25').Value = '  +0.28%  '
$ws.Range('D26').Value = '155.02'
$ws.Range('E26').Value = '  -0.16%  '
$ws.Range('D27').Value = '18.69'
$ws.Range('E27').Value = '  +0.92%  '
$ws.Range('D28').Value = '5.148'
$ws.Range('E28').Value = '  -1.44%  '
$ws.Range('D29').Value = '114.21'
$ws.Range('E29').Value = '  -4.92%  '
$ws.Range('D30').Value = '1.797'
$ws.Range('E30').Value = '  -4.08%  '
$ws.Range('D31').Value = '0.08853'
$ws.Range('E31').Value = '  -0.52%  '
$ws.Range('D32').Value = '0.7467'
$ws.Range('E32').Value = '  -2.31%  '
$ws.Range('D33').Value = '2.977'
$ws.Range('E33').Value = '  +0.47%  '
$ws.Range('D34').Value = '4.534'
$ws.Range('E34').Value = '  +0.82%  '
$ws.Range('D35').Value = '1.119'
$ws.Range('E35').Value = '  -0.51%  '
$ws.Range('D36').Value = '1.001'
$ws.Range('E36').Value = '  -0.04%  '
$ws.Range('D37').Value = '1.086'
$ws.Range('E37').Value = '  -1.58%  '
$ws.Range('D38').Value = '0.05285'
$ws.Range('E38').Value = '  -2.64%  '
$ws.Range('D39').Value = '0.01918'
$ws.Range('E39').Value = '  -0.42%  '
$ws.Range('E40').Value = '  -1.33%  '
$ws.Range('D41').Value = '0.5064'
$ws.Range('E41').Value = '  -0.04%  '
$ws.Range('D42').Value = '0.1642'
$ws.Range('E42').Value = '  -1.14%  '
$ws.Range('D43').Value = '6.436'
$ws.Range('E43').Value = '  -1.94%  '
$ws.Range('D44').Value = '8.343'
$ws.Range('E44').Value = '  -0.67%  '
$ws.Range('D45').Value = '10.39'
$ws.Range('E45').Value = '  +0.20%  '
$ws.Range('D46').Value = '105.62'
$ws.Range('E46').Value = '  -0.55%  '
$ws.Range('D47').Value = '0.4676'
$ws.Range('E47').Value = '  +0.85%  '
$ws.Range('D48').Value = '0.06445'
$ws.Range('E48').Value = '  -1.54%  '
$ws.Range('E49').Value = '  -0.11%  '
$ws.Range('E50').Value = '  -1.27%  '
$ws.Range('D51').Value = '63.56'
$ws.Range('E51').Value = '  -1.30%  '
